# forests-scraped.xlsx update
#
# "New" sheet rows 2-9 (the 8 previously-new listings) graduate into the
# "Previously added" sheet (appended after the existing last row), and the
# "New" sheet is repopulated with 3 freshly scraped listings.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Previously added")
$ws2 = $wb.Worksheets.Item("New")

# ---------------------------------------------------------------------
# 1) Move the 8 rows currently on "New" (rows 2-9) down onto the bottom
#    of "Previously added" (rows 273-280). A straight range copy carries
#    over both the values (kept as text / shared-strings) and the cell
#    styles in one shot.
# ---------------------------------------------------------------------
$srcRows = $ws2.Range("A2:F9")
$dstRows = $ws1.Range("A273:F280")
$srcRows.Copy($dstRows)
$excel.CutCopyMode = $false

# Hyperlinks don't travel with a cross-sheet Range.Copy, so re-create them
# on the destination rows - the visible text of column A *is* the URL.
for ($r = 273; $r -le 280; $r++) {
    $addr = $ws1.Cells.Item($r, 1).Value()
    $ws1.Hyperlinks.Add($ws1.Cells.Item($r, 1), $addr)
}

# ---------------------------------------------------------------------
# 2) Clear out the old "New" sheet content.
#    (Range.Hyperlinks.Delete() drops every hyperlink on the sheet, so do
#    it before we rebuild rows 2-4 with their own fresh links.)
# ---------------------------------------------------------------------
$ws2.Range("A2").Hyperlinks.Delete()

# Drop rows 5-9 entirely (shifts the used range / dimension back to F4).
$ws2.Rows("5:9").Delete()

# ---------------------------------------------------------------------
# 3) Write the 3 new listings into "New" rows 2-4.
#    Helper below forces a numeric-looking cadastre number to stay text
#    (matching the source data's shared-string storage) without touching
#    the cell's existing style: write it as a quoted-string formula, then
#    paste-special just the computed value back over itself.
# ---------------------------------------------------------------------
function Set-TextValue($cell, [string]$text) {
    $escaped = $text.Replace('"', '""')
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
    $excel.CutCopyMode = $false
}

# Row 2
$ws2.Range("A2").Value = "https://www.ss.com/msg/lv/real-estate/wood/balvi-and-reg/balvu-pag/akimk.html"
$ws2.Range("B2").Value = "15 000 €"
$ws2.Range("C2").Value = "Balvi un raj."
$ws2.Range("D2").Value = "1 ha."
Set-TextValue $ws2.Range("E2") "38460040026"
$ws2.Range("F2").Value = 45985.790972222225

# Row 3
$ws2.Range("A3").Value = "https://www.ss.com/msg/lv/real-estate/wood/dobele-and-reg/dobeles-pag/lfnll.html"
$ws2.Range("B3").Value = "35 000 €"
$ws2.Range("C3").Value = "Dobele un raj."
$ws2.Range("D3").Value = "6 ha."
Set-TextValue $ws2.Range("E3") "46600010112"
$ws2.Range("F3").Value = 45985.72777777778

# Row 4
$ws2.Range("A4").Value = "https://www.ss.com/msg/lv/real-estate/wood/madona-and-reg/aronas-pag/cghlhb.html"
$ws2.Range("B4").Value = "80 000 €"
$ws2.Range("C4").Value = "Madona un raj."
$ws2.Range("D4").Value = "137000 m²"
$ws2.Range("E4").Value = ""
$ws2.Range("F4").Value = 45985.65902777778

# Re-add hyperlinks for the 3 new rows.
$ws2.Hyperlinks.Add($ws2.Cells.Item(2, 1), $ws2.Cells.Item(2, 1).Value())
$ws2.Hyperlinks.Add($ws2.Cells.Item(3, 1), $ws2.Cells.Item(3, 1).Value())
$ws2.Hyperlinks.Add($ws2.Cells.Item(4, 1), $ws2.Cells.Item(4, 1).Value())
